# Auto-generated Excel COM-interop edit script
# Applies the weekly crime-data refresh described in the commit diff.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header text updates (report volume/number and week-covering dates) ---
$ws.Range("A8").Value = "Volume 30   Number  43"
$ws.Range("C9").Value = "Report Covering the Week  10/23/2023  Through  10/29/2023"

# --- Cells that change from a number to the text placeholder "0" or "***.*" ---
# Donor cells C28 (text "0", style 14) and E28 (text "***.*", style 14) are untouched
# elsewhere in this edit, so copying their format+value over gives an exact match.
$ws.Range("C28").Copy($ws.Range("C15"))
$ws.Range("C28").Copy($ws.Range("C18"))
$ws.Range("C28").Copy($ws.Range("C22"))
$ws.Range("C28").Copy($ws.Range("C26"))
$ws.Range("C28").Copy($ws.Range("D27"))
$ws.Range("E28").Copy($ws.Range("E27"))
$ws.Range("C28").Copy($ws.Range("F30"))

# --- Cells that change from the text placeholder to a real number ---
# Donor cells I28 (style 15) and K28 (style 16) are untouched elsewhere, so copying
# their format gives the right numeric style; the real value is set right after.
$ws.Range("I28").Copy($ws.Range("D20"))
$ws.Range("D20").Value = 3
$ws.Range("K28").Copy($ws.Range("E20"))
$ws.Range("E20").Value = 66.666666666666
$ws.Range("I28").Copy($ws.Range("D22"))
$ws.Range("D22").Value = 1
$ws.Range("K28").Copy($ws.Range("E22"))
$ws.Range("E22").Value = -100
$ws.Range("I28").Copy($ws.Range("D30"))
$ws.Range("D30").Value = 1
$ws.Range("K28").Copy($ws.Range("E30"))
$ws.Range("E30").Value = -100
$ws.Range("I28").Copy($ws.Range("G30"))
$ws.Range("G30").Value = 1
$ws.Range("K28").Copy($ws.Range("H30"))
$ws.Range("H30").Value = -100

# --- Plain value updates (style/type unchanged) ---
$ws.Range("G15").Value = 1
$ws.Range("H15").Value = 0
$ws.Range("C16").Value = 1
$ws.Range("E16").Value = -50
$ws.Range("G16").Value = 8
$ws.Range("H16").Value = 50
$ws.Range("I16").Value = 123
$ws.Range("J16").Value = 89
$ws.Range("K16").Value = 38.202247191011
$ws.Range("L16").Value = 105
$ws.Range("M16").Value = -10.218978102189
$ws.Range("N16").Value = -80.959752321981
$ws.Range("C17").Value = 7
$ws.Range("E17").Value = 600
$ws.Range("F17").Value = 18
$ws.Range("G17").Value = 7
$ws.Range("H17").Value = 157.142857142857
$ws.Range("I17").Value = 156
$ws.Range("J17").Value = 124
$ws.Range("K17").Value = 25.806451612903
$ws.Range("L17").Value = 14.705882352941
$ws.Range("M17").Value = 64.210526315789
$ws.Range("N17").Value = -36.326530612244
$ws.Range("D18").Value = 3
$ws.Range("E18").Value = -100
$ws.Range("J18").Value = 103
$ws.Range("K18").Value = -28.155339805825
$ws.Range("L18").Value = 37.037037037037
$ws.Range("N18").Value = -85.258964143426
$ws.Range("C19").Value = 8
$ws.Range("D19").Value = 5
$ws.Range("E19").Value = 60
$ws.Range("F19").Value = 25
$ws.Range("G19").Value = 25
$ws.Range("I19").Value = 327
$ws.Range("J19").Value = 285
$ws.Range("K19").Value = 14.736842105263
$ws.Range("L19").Value = 50
$ws.Range("M19").Value = 43.421052631578
$ws.Range("N19").Value = -41.814946619217
$ws.Range("C20").Value = 5
$ws.Range("F20").Value = 15
$ws.Range("G20").Value = 10
$ws.Range("H20").Value = 50
$ws.Range("I20").Value = 97
$ws.Range("J20").Value = 66
$ws.Range("K20").Value = 46.969696969697
$ws.Range("L20").Value = 125.581395348837
$ws.Range("M20").Value = 340.909090909091
$ws.Range("N20").Value = -69.303797468354
$ws.Range("C21").Value = 21
$ws.Range("D21").Value = 14
$ws.Range("E21").Value = 50
$ws.Range("F21").Value = 74
$ws.Range("G21").Value = 58
$ws.Range("H21").Value = 27.586206896551
$ws.Range("I21").Value = 785
$ws.Range("J21").Value = 683
$ws.Range("K21").Value = 14.934114202049
$ws.Range("L21").Value = 52.131782945736
$ws.Range("M21").Value = 40.178571428571
$ws.Range("N21").Value = -65.928819444444
$ws.Range("F22").Value = 3
$ws.Range("G22").Value = 3
$ws.Range("H22").Value = 0
$ws.Range("J22").Value = 17
$ws.Range("K22").Value = 88.235294117647
$ws.Range("M22").Value = 18.518518518518
$ws.Range("C23").Value = 3
$ws.Range("E23").Value = 200
$ws.Range("F23").Value = 14
$ws.Range("G23").Value = 6
$ws.Range("H23").Value = 133.333333333333
$ws.Range("I23").Value = 129
$ws.Range("J23").Value = 96
$ws.Range("K23").Value = 34.375
$ws.Range("L23").Value = 41.758241758241
$ws.Range("M23").Value = 89.705882352941
$ws.Range("C24").Value = 10
$ws.Range("D24").Value = 8
$ws.Range("E24").Value = 25
$ws.Range("F24").Value = 40
$ws.Range("G24").Value = 41
$ws.Range("H24").Value = -2.439024390243
$ws.Range("I24").Value = 429
$ws.Range("J24").Value = 429
$ws.Range("K24").Value = 0
$ws.Range("L24").Value = 21.186440677966
$ws.Range("M24").Value = -18.285714285714
$ws.Range("C25").Value = 9
$ws.Range("E25").Value = 200
$ws.Range("F25").Value = 22
$ws.Range("G25").Value = 13
$ws.Range("H25").Value = 69.230769230769
$ws.Range("I25").Value = 204
$ws.Range("J25").Value = 177
$ws.Range("K25").Value = 15.254237288135
$ws.Range("L25").Value = 17.241379310344
$ws.Range("M25").Value = -18.725099601593
$ws.Range("G26").Value = 2
$ws.Range("H26").Value = 0
$ws.Range("F27").Value = 3
$ws.Range("G27").Value = 11
$ws.Range("H27").Value = -72.727272727272
$ws.Range("L27").Value = 14.285714285714
$ws.Range("J30").Value = 10
$ws.Range("K30").Value = -50
